# Insert a new row at position 22 (pushes existing rows 22-55 down to 23-56)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with its data.
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value = "Bíobío"
$ws.Cells.Item(22, 4).Value = 44540
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = 100112001
$ws.Cells.Item(22, 7).Value = "Berenjena"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 110
$ws.Cells.Item(22, 11).Value = 9000
$ws.Cells.Item(22, 12).Value = 10000
$ws.Cells.Item(22, 13).Value = 9455
$ws.Cells.Item(22, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 158
$ws.Cells.Item(22, 17).Value = 60
$ws.Cells.Item(22, 18).Value = "Hortaliza"
